# Updates cryptos list data (prices + 1h volume %) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.904.40'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').Value = '1.643.43'
$ws.Range('E3').Value = '  +1.23%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'213.55"
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'23.52"
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('E10').Value = '  +0.70%  '
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('D12').Value = '1.876.31'
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('D13').Value = '1.646.44'
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('E14').Value = '  +4.08%  '
$ws.Range('E15').Value = '  +0.51%  '
$ws.Range('D16').Value = "'65.87"
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('D17').Value = '27.888.03'
$ws.Range('E17').Value = '  +1.31%  '
$ws.Range('D18').Value = "'230.51"
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('D19').Value = '0.0₃0725'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').Value = "'7.63"
$ws.Range('E20').Value = '  +0.72%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').Value = "'10.83"
$ws.Range('E22').Value = '  +5.27%  '
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('E24').Value = '  +2.53%  '
$ws.Range('D25').Value = "'152.74"
$ws.Range('E25').Value = '  +1.83%  '
$ws.Range('E26').Value = '  +0.58%  '
$ws.Range('E27').Value = '  +1.07%  '
$ws.Range('D28').Value = "'15.72"
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('E30').Value = '  +1.11%  '
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('E32').Value = '  +1.88%  '
$ws.Range('D33').Value = '1.433.75'
$ws.Range('E33').Value = '  -2.74%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('E35').Value = '  +1.68%  '
$ws.Range('D36').Value = "'2.34"
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('E37').Value = '  +1.54%  '
$ws.Range('D38').Value = "'0.0168"
$ws.Range('E38').Value = '  +0.77%  '
$ws.Range('D39').Value = "'0.929"
$ws.Range('E39').Value = '  -2.63%  '
$ws.Range('D40').Value = "'0.558"
$ws.Range('E40').Value = '  +0.89%  '
$ws.Range('E41').Value = '  +1.95%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').Value = "'68.59"
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('E44').Value = '  +0.47%  '
$ws.Range('D45').Value = "'5.44"
$ws.Range('E45').Value = '  +3.25%  '
$ws.Range('E46').Value = '  +3.52%  '
$ws.Range('E47').Value = '  +0.16%  '
$ws.Range('D48').Value = '1.785.22'
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('D49').Value = "'89.17"
$ws.Range('E49').Value = '  +2.09%  '
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.0506"
$ws.Range('E51').Value = '  +0.55%  '
